$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new travel city rows (city, country, lat, lon) starting at row 89
$newCities = @(
    @("Baltimore", "USA", 39.2904, -76.6122),
    @("Yachats", "USA", 44.3112, -124.1048),
    @("Newport", "USA", 44.6368, -125.0535),
    @("Joshua Tree National Park", "USA", 33.8734, -115.901),
    @("Bend", "USA", 44.0582, -121.3153)
)

$row = 89
foreach ($city in $newCities) {
    $ws.Cells.Item($row, 1).Value = $city[0]
    $ws.Cells.Item($row, 2).Value = $city[1]
    $ws.Cells.Item($row, 3).Value = $city[2]
    $ws.Cells.Item($row, 4).Value = $city[3]
    $row++
}

# Widen column A slightly to fit the longer city names
$ws.Columns.Item(1).ColumnWidth = 20.8

# Update the view: scroll back to top and select A12 (matches the saved cursor position)
$ws.Range("A12").Select() | Out-Null
